# Apply the crypto price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel auto-parses pure numeric text back into a Number, which would drop
# formatting like trailing zeros ('1.00' -> 1). The sheet stores Price as
# text, so force-prefix an apostrophe for values that parse as plain numbers
# to keep them stored as text, matching the source data.
$apos = "'"

# Row 2
$ws.Range("D2").Value = '64.332.94'
$ws.Range("E2").Value = '  -0.06%  '

# Row 3
$ws.Range("D3").Value = '3.143.52'
$ws.Range("E3").Value = '  -0.65%  '

# Row 4
$ws.Range("D4").Value = $apos + '0.999'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").Value = $apos + '610.29'
$ws.Range("E5").Value = '  +0.59%  '

# Row 6
$ws.Range("D6").Value = $apos + '143.20'
$ws.Range("E6").Value = '  -3.50%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = '3.141.05'
$ws.Range("E8").Value = '  -0.72%  '

# Row 9
$ws.Range("E9").Value = '  -1.57%  '

# Row 10
$ws.Range("D10").Value = $apos + '0.149'
$ws.Range("E10").Value = '  -2.17%  '

# Row 11
$ws.Range("D11").Value = $apos + '5.32'
$ws.Range("E11").Value = '  -4.69%  '

# Row 12
$ws.Range("E12").Value = '  -2.52%  '

# Row 13
$ws.Range("E13").Value = '  -1.54%  '

# Row 14
$ws.Range("D14").Value = $apos + '35.12'
$ws.Range("E14").Value = '  -3.81%  '

# Row 15
$ws.Range("D15").Value = '3.654.66'
$ws.Range("E15").Value = '  -0.91%  '

# Row 16
$ws.Range("E16").Value = '  +2.44%  '

# Row 17
$ws.Range("D17").Value = '64.296.21'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18
$ws.Range("D18").Value = '3.134.97'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
$ws.Range("E19").Value = '  -2.57%  '

# Row 20
$ws.Range("D20").Value = $apos + '471.14'
$ws.Range("E20").Value = '  -2.37%  '

# Row 21
$ws.Range("D21").Value = $apos + '14.48'
$ws.Range("E21").Value = '  -1.01%  '

# Row 22
$ws.Range("D22").Value = $apos + '0.716'
$ws.Range("E22").Value = '  +0.08%  '

# Row 23
$ws.Range("D23").Value = $apos + '7.78'
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
$ws.Range("D24").Value = $apos + '13.68'
$ws.Range("E24").Value = '  -0.97%  '

# Row 25
$ws.Range("D25").Value = $apos + '83.52'
$ws.Range("E25").Value = '  -0.15%  '

# Row 26
$ws.Range("D26").Value = $apos + '1.00'
$ws.Range("E26").Value = '  +0.13%  '

# Row 27
$ws.Range("D27").Value = $apos + '2.78'
$ws.Range("E27").Value = '  -4.45%  '

# Row 28
$ws.Range("D28").Value = $apos + '8.42'
$ws.Range("E28").Value = '  -1.41%  '

# Row 29
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").Value = $apos + '7.19'
$ws.Range("E29").Value = '  +3.58%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = $apos + '2.08'
$ws.Range("E30").Value = '  -6.66%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = $apos + '0.115'
$ws.Range("E31").Value = '  -5.26%  '

# Row 32
$ws.Range("E32").Value = '  -0.18%  '

# Row 33
$ws.Range("D33").Value = $apos + '26.19'
$ws.Range("E33").Value = '  -0.90%  '

# Row 34
$ws.Range("D34").Value = $apos + '2.61'
$ws.Range("E34").Value = '  -6.60%  '

# Row 35
$ws.Range("E35").Value = '  +0.65%  '

# Row 36
$ws.Range("D36").Value = $apos + '5.90'
$ws.Range("E36").Value = '  -3.07%  '

# Row 37
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = $apos + '52.79'
$ws.Range("E37").Value = '  -2.99%  '

# Row 38
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0750'
$ws.Range("E38").Value = '  +1.73%  '

# Row 39
$ws.Range("D39").Value = $apos + '3.04'
$ws.Range("E39").Value = '  -0.41%  '

# Row 40
$ws.Range("D40").Value = $apos + '449.66'
$ws.Range("E40").Value = '  -1.56%  '

# Row 41
$ws.Range("D41").Value = $apos + '0.0391'
$ws.Range("E41").Value = '  -2.39%  '

# Row 42
$ws.Range("D42").Value = $apos + '0.117'
$ws.Range("E42").Value = '  -5.31%  '

# Row 43
$ws.Range("D43").Value = $apos + '8.24'
$ws.Range("E43").Value = '  -2.94%  '

# Row 44
$ws.Range("D44").Value = '2.822.18'
$ws.Range("E44").Value = '  -2.16%  '

# Row 45
$ws.Range("D45").Value = $apos + '2.26'
$ws.Range("E45").Value = '  -1.58%  '

# Row 46
$ws.Range("D46").Value = $apos + '0.262'
$ws.Range("E46").Value = '  -4.00%  '

# Row 47
$ws.Range("D47").Value = $apos + '2.41'
$ws.Range("E47").Value = '  +2.75%  '

# Row 48
$ws.Range("D48").Value = $apos + '0.999'
$ws.Range("E48").Value = '  +0.08%  '

# Row 49
$ws.Range("D49").Value = $apos + '26.18'
$ws.Range("E49").Value = '  -1.45%  '

# Row 50
$ws.Range("E50").Value = '  -2.12%  '

# Row 51
$ws.Range("D51").Value = $apos + '34.43'
$ws.Range("E51").Value = '  +3.55%  '
